# Auto-generated edit script
# Applies the cell-value changes described by the commit diff to the
# per-job "Diabolos_Profits" market-data tables (one table per crafting job sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 852.7273
$ws.Range("J12").Value = 497.5
$ws.Range("L12").Value = 497.5
$ws.Range("N12").Value = -837.5
$ws.Range("H74").Value = 4542.857
$ws.Range("J74").Value = 4600
$ws.Range("L74").Value = 4600
$ws.Range("N74").Value = -6472
$ws.Range("H77").Value = 4542.857
$ws.Range("J77").Value = 4600
$ws.Range("L77").Value = 23000
$ws.Range("N77").Value = -32360
$ws.Range("H100").Value = 5725.5
$ws.Range("I100").Value = 5695
$ws.Range("K100").Value = 5695
$ws.Range("M100").Value = -5154
$ws.Range("H111").Value = 49027.43
$ws.Range("I111").Value = 27032.25
$ws.Range("J111").Value = 78354.336
$ws.Range("K111").Value = 81096.75
$ws.Range("L111").Value = 235063.008
$ws.Range("M111").Value = -78029.75
$ws.Range("N111").Value = -241197.008
$ws.Range("H113").Value = 90913260
$ws.Range("I113").Value = 333337470
$ws.Range("J113").Value = 4186.375
$ws.Range("K113").Value = 333337470
$ws.Range("L113").Value = 4186.375
$ws.Range("M113").Value = -333334216
$ws.Range("N113").Value = -10694.375
$ws.Range("H116").Value = 38969684
$ws.Range("J116").Value = 55563900
$ws.Range("L116").Value = 55563900
$ws.Range("N116").Value = -55570784
$ws.Range("H132").Value = 3281.5286
$ws.Range("I132").Value = 3076.4897
$ws.Range("J132").Value = 3759.9524
$ws.Range("K132").Value = 9229.4691
$ws.Range("L132").Value = 11279.8572
$ws.Range("M132").Value = -6699.4691
$ws.Range("N132").Value = -16339.8572
$ws.Range("H137").Value = 55560396
$ws.Range("J137").Value = 6574.3
$ws.Range("L137").Value = 19722.9
$ws.Range("N137").Value = -24822.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3236.9666
$ws.Range("I74").Value = 3407.6155
$ws.Range("K74").Value = 3407.6155
$ws.Range("M74").Value = -2533.6155
$ws.Range("H77").Value = 3236.9666
$ws.Range("I77").Value = 3407.6155
$ws.Range("K77").Value = 17038.0775
$ws.Range("M77").Value = -12670.0775
$ws.Range("H110").Value = 166692880
$ws.Range("I110").Value = 200001440
$ws.Range("K110").Value = 200001440
$ws.Range("M110").Value = -199999395

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22729250
$ws.Range("I86").Value = 33334900
$ws.Range("J86").Value = 2858.1428
$ws.Range("K86").Value = 33334900
$ws.Range("L86").Value = 2858.1428
$ws.Range("M86").Value = -33333777
$ws.Range("N86").Value = -5104.1428
$ws.Range("H89").Value = 22729250
$ws.Range("I89").Value = 33334900
$ws.Range("J89").Value = 2858.1428
$ws.Range("K89").Value = 166674500
$ws.Range("L89").Value = 14290.714
$ws.Range("M89").Value = -166668884
$ws.Range("N89").Value = -25522.714
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H134").Value = 2160
$ws.Range("I134").Value = 1700.8334
$ws.Range("K134").Value = 5102.5002
$ws.Range("M134").Value = -2567.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2988
$ws.Range("I16").Value = 2988
$ws.Range("K16").Value = 2988
$ws.Range("M16").Value = -2701
$ws.Range("H99").Value = 2933.7693
$ws.Range("I99").Value = 2602.4
$ws.Range("J99").Value = 3140.875
$ws.Range("K99").Value = 2602.4
$ws.Range("L99").Value = 3140.875
$ws.Range("M99").Value = -1104.4
$ws.Range("N99").Value = -6136.875
$ws.Range("H113").Value = 2988
$ws.Range("I113").Value = 2988
$ws.Range("K113").Value = 2988
$ws.Range("M113").Value = -818
$ws.Range("H122").Value = 6449.75
$ws.Range("I122").Value = 2900
$ws.Range("J122").Value = 9999.5
$ws.Range("K122").Value = 8700
$ws.Range("L122").Value = 29998.5
$ws.Range("M122").Value = -6250
$ws.Range("N122").Value = -34898.5
$ws.Range("H126").Value = 2933.7693
$ws.Range("I126").Value = 2602.4
$ws.Range("J126").Value = 3140.875
$ws.Range("K126").Value = 7807.200000000001
$ws.Range("L126").Value = 9422.625
$ws.Range("M126").Value = -5337.200000000001
$ws.Range("N126").Value = -14362.625
$ws.Range("H132").Value = 4370.357
$ws.Range("J132").Value = 4828
$ws.Range("L132").Value = 14484
$ws.Range("N132").Value = -19544
$ws.Range("H141").Value = 79582.336
$ws.Range("J141").Value = 79582.336
$ws.Range("L141").Value = 79582.336
$ws.Range("N141").Value = -89942.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 30732.684
$ws.Range("I120").Value = 17944.75
$ws.Range("J120").Value = 40033
$ws.Range("K120").Value = 53834.25
$ws.Range("L120").Value = 120099
$ws.Range("M120").Value = -48996.25
$ws.Range("N120").Value = -129775
$ws.Range("H122").Value = 382
$ws.Range("J122").Value = 500
$ws.Range("L122").Value = 4500
$ws.Range("N122").Value = -9400
$ws.Range("H129").Value = 2887.625
$ws.Range("J129").Value = 2887.625
$ws.Range("L129").Value = 8662.875
$ws.Range("N129").Value = -18662.875
$ws.Range("H137").Value = 1141.7778
$ws.Range("I137").Value = 909.5
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 2728.5
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = 2371.5
$ws.Range("N137").Value = -19200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3479.8948
$ws.Range("I80").Value = 2694.6667
$ws.Range("K80").Value = 2694.6667
$ws.Range("M80").Value = -1696.6667
$ws.Range("H83").Value = 3479.8948
$ws.Range("I83").Value = 2694.6667
$ws.Range("K83").Value = 13473.3335
$ws.Range("M83").Value = -8481.333500000001
$ws.Range("H103").Value = 52000
$ws.Range("J103").Value = 52000
$ws.Range("L103").Value = 52000
$ws.Range("N103").Value = -54344
$ws.Range("H107").Value = 790
$ws.Range("I107").Value = 601.4286
$ws.Range("K107").Value = 601.4286
$ws.Range("M107").Value = 1318.5714
$ws.Range("H113").Value = 2856.423
$ws.Range("I113").Value = 1761.7778
$ws.Range("J113").Value = 3435.9412
$ws.Range("K113").Value = 1761.7778
$ws.Range("L113").Value = 3435.9412
$ws.Range("M113").Value = 408.2221999999999
$ws.Range("N113").Value = -7775.9412
$ws.Range("H122").Value = 2725.8635
$ws.Range("J122").Value = 3976.3333
$ws.Range("L122").Value = 11928.9999
$ws.Range("N122").Value = -16828.9999
$ws.Range("H127").Value = 90313
$ws.Range("J127").Value = 90313
$ws.Range("L127").Value = 90313
$ws.Range("N127").Value = -100233
$ws.Range("H132").Value = 404823.1
$ws.Range("I132").Value = 718112.5
$ws.Range("K132").Value = 2154337.5
$ws.Range("M132").Value = -2151807.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 945.7143
$ws.Range("I22").Value = 915
$ws.Range("J22").Value = 954.0909
$ws.Range("K22").Value = 915
$ws.Range("L22").Value = 954.0909
$ws.Range("M22").Value = -620
$ws.Range("N22").Value = -1544.0909
$ws.Range("H27").Value = 945.7143
$ws.Range("I27").Value = 915
$ws.Range("J27").Value = 954.0909
$ws.Range("K27").Value = 915
$ws.Range("L27").Value = 954.0909
$ws.Range("M27").Value = -808
$ws.Range("N27").Value = -1168.0909
$ws.Range("H40").Value = 3699.9333
$ws.Range("I40").Value = 2944.3333
$ws.Range("J40").Value = 4833.3335
$ws.Range("K40").Value = 2944.3333
$ws.Range("L40").Value = 4833.3335
$ws.Range("M40").Value = -2808.3333
$ws.Range("N40").Value = -5105.3335
$ws.Range("H124").Value = 79991
$ws.Range("J124").Value = 79991
$ws.Range("L124").Value = 79991
$ws.Range("N124").Value = -89811
$ws.Range("H136").Value = 7425.0625
$ws.Range("I136").Value = 4207.5557
$ws.Range("J136").Value = 11561.857
$ws.Range("K136").Value = 12622.6671
$ws.Range("L136").Value = 34685.571
$ws.Range("M136").Value = -10072.6671
$ws.Range("N136").Value = -39785.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3929.6
$ws.Range("I96").Value = 3562
$ws.Range("J96").Value = 5400
$ws.Range("K96").Value = 3562
$ws.Range("L96").Value = 5400
$ws.Range("M96").Value = -2189
$ws.Range("N96").Value = -8146
$ws.Range("H100").Value = 399
$ws.Range("I100").Value = 399
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 798
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -257
$ws.Range("N100").ClearContents()
$ws.Range("H129").Value = 69999.664
$ws.Range("I129").Value = 79999.5
$ws.Range("J129").Value = 50000
$ws.Range("K129").Value = 79999.5
$ws.Range("L129").Value = 50000
$ws.Range("M129").Value = -74999.5
$ws.Range("N129").Value = -60000
$ws.Range("H132").Value = 292661.03
$ws.Range("I132").Value = 387062.44
$ws.Range("J132").Value = 3903.7646
$ws.Range("K132").Value = 1161187.32
$ws.Range("L132").Value = 11711.2938
$ws.Range("M132").Value = -1158657.32
$ws.Range("N132").Value = -16771.2938
$ws.Range("H135").Value = 99999
$ws.Range("I135").Value = 99999
$ws.Range("K135").Value = 99999
$ws.Range("M135").Value = -94929
